# Auto update stock data
# Update the "Date_1" column (column A) from 2026/01/03 to 2026/01/04
# for each company's most recent data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2026/01/03") {
        # Temporarily mark the cell as Text so Excel doesn't
        # auto-convert the date-looking string into a date serial
        # number, then restore the original (unstyled/"Normal") look.
        $cell.NumberFormat = "@"
        $cell.Value = "2026/01/04"
        $cell.Style = "Normal"
    }
}
